$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.947.77'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '1.714.61'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.99'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3968'
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4116'
$ws.Range("E8").Value = '  +1.86%  '
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.004'
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.36'
$ws.Range("E11").Value = '  +4.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08934'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.732'
$ws.Range("E13").Value = '  +6.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.00'
$ws.Range("E14").Value = '  +6.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.153'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001370'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").Value = '1.697.18'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.55'
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07148'
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.508'
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.52'
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("D24").Value = '24.955.35'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.140'
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.340'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.25'
$ws.Range("E27").Value = '  +1.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.235'
$ws.Range("E28").Value = '  +23.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.21'
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '140.98'
$ws.Range("E30").Value = '  +2.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.241'
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.885'
$ws.Range("E32").Value = '  +10.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09018'
$ws.Range("E33").Value = '  +4.44%  '
$ws.Range("D34").Value = '1.880.79'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.087'
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03006'
$ws.Range("E36").Value = '  +10.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2807'
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.17'
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.968'
$ws.Range("E39").Value = '  +2.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '14.59'
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09316'
$ws.Range("E41").Value = '  +1.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8142'
$ws.Range("E42").Value = '  +5.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.488'
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.78'
$ws.Range("E44").Value = '  +6.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7374'
$ws.Range("E45").Value = '  +2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.650'
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.272'
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.348'
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.90'
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '93.68'
$ws.Range("E51").Value = '  +4.28%  '
